# Update the crypto price/volume table (rows 2-51) on Sheet1 to the latest
# scraped values. Numeric-looking price strings are written with a leading
# apostrophe so Excel keeps them as text (matching the workbook's original
# inlineStr cell type) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.111.52"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.471.54"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'583.63"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'131.18"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "'7.60"
$ws.Range("E9").Value = "  +4.91%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "4.067.71"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "3.474.42"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "64.088.08"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'24.23"
$ws.Range("E17").Value = "  -6.85%  "
$ws.Range("D18").Value = "'9.95"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'5.67"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "'13.40"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'383.62"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "'0.567"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "3.614.11"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "'74.97"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "'7.05"
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("D33").Value = "3.501.48"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D36").Value = "'22.85"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").Value = "'5.17"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "'6.73"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'162.07"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.49"
$ws.Range("E40").Value = "  -4.14%  "
$ws.Range("D41").Value = "'0.0776"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'41.30"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").Value = "'1.61"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").Value = "'23.36"
$ws.Range("E48").Value = "  -7.50%  "
$ws.Range("D49").Value = "'6.69"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").Value = "'0.901"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "2.326.32"
$ws.Range("E51").Value = "  -5.61%  "
